# Updated cryptos list on Mon Feb 20 22:44:09 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row, and
# for a handful of rows whose rank order changed this run, also rewrites
# the Coin (B) / Link (C) columns so the row reflects the new coin.
#
# Numeric-looking Price strings (e.g. "1.002", "0.08380") are written with
# a leading apostrophe so Excel stores them as literal text (preserving
# exact formatting / trailing zeros) instead of silently parsing them into
# floating point numbers; the Style is then reset to "Normal" so no stray
# number-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # NOTE: positional params only - this COM-interop PowerShell host does
    # not reliably bind named (-Cell / -Value) arguments on user functions.
    param(
        [string]$Cell,
        [string]$Value
    )
    $range = $ws.Range($Cell)
    # A Price string that parses as a plain decimal number (e.g. "1.002",
    # "0.08380") would otherwise get silently coerced into a float by
    # Excel on assignment (losing trailing zeros / exact text). Force it
    # to stay text via the classic leading-apostrophe trick, then strip
    # the resulting "quote prefix" style so no stray style index sticks
    # around on the cell.
    if ($Value -match '^[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $Value
        $range.Style = "Normal"
    } else {
        $range.Value = $Value
    }
}

# row -> @{ B=..; C=..; D=..; E=.. }  (only keys that changed are present)
$rowUpdates = [ordered]@{
    2  = @{ D = "24.740.56";  E = "  +0.59%  " }
    3  = @{ D = "1.699.16";   E = "  +0.54%  " }
    4  = @{ D = "1.002";      E = "  -0.12%  " }
    5  = @{ D = "314.55";     E = "  +0.08%  " }
    6  = @{ E = "  -0.15%  " }
    7  = @{ D = "0.3979";     E = "  +2.33%  " }
    8  = @{ D = "0.4049";     E = "  +0.38%  " }
    9  = @{ E = "  +0.04%  " }
    10 = @{ D = "53.57";      E = "  +0.88%  " }
    11 = @{ D = "1.463";      E = "  -2.02%  " }
    12 = @{ D = "0.08796";    E = "  +0.61%  " }
    13 = @{ D = "26.21";      E = "  +3.26%  " }
    14 = @{ D = "7.516";      E = "  -0.08%  " }
    15 = @{ D = "7.961";      E = "  +0.22%  " }
    16 = @{ D = "0.00001339"; E = "  -1.29%  " }
    17 = @{ D = "1.689.65";   E = "  -0.16%  " }
    18 = @{ D = "95.48";      E = "  -3.07%  " }
    19 = @{ D = "0.07173";    E = "  +0.96%  " }
    20 = @{ D = "20.82";      E = "  +4.42%  " }
    21 = @{ D = "7.317";      E = "  +0.80%  " }
    22 = @{ D = "1.003";      E = "  +0.05%  " }
    23 = @{ D = "14.36";      E = "  +0.83%  " }
    24 = @{ D = "24.729.33";  E = "  +0.57%  " }
    25 = @{ D = "2.382";      E = "  +1.15%  " }
    26 = @{ E = "  -3.35%  " }
    27 = @{ D = "23.09";      E = "  +1.63%  " }
    28 = @{ D = "6.106";      E = "  +16.98%  " }
    29 = @{ D = "162.03";     E = "  +0.10%  " }
    30 = @{ D = "143.90";     E = "  +5.16%  " }
    31 = @{ D = "8.265";      E = "  -5.85%  " }

    32 = @{
        B = "WEMIXTOKEN"
        C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
        D = "2.260"
        E = "  +15.29%  "
    }
    33 = @{
        B = "WrappedliquidstakedEther2.0"
        C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
        D = "1.913.35"
        E = "  +1.97%  "
    }
    34 = @{
        B = "Hedera"
        C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
        D = "0.08575"
        E = "  -2.89%  "
    }
    35 = @{
        B = "InternetComputer(DFINITY)"
        C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
        D = "7.306"
        E = "  -0.97%  "
    }
    36 = @{
        B = "VeChain"
        C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
        D = "0.03171"
        E = "  +8.93%  "
    }

    37 = @{ D = "1.027";   E = "  -0.63%  " }
    38 = @{ D = "0.2843";  E = "  +3.46%  " }

    39 = @{
        B = "Stellar"
        C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
        D = "0.09438"
        E = "  +3.44%  "
    }
    40 = @{
        B = "TheSandbox"
        C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
        D = "0.8288"
        E = "  +4.88%  "
    }

    41 = @{ D = "10.68";    E = "  -0.82%  " }
    42 = @{ D = "14.15";    E = "  -0.53%  " }
    43 = @{ D = "1.477";    E = "  +1.33%  " }
    44 = @{ D = "17.57";    E = "  +4.87%  " }
    45 = @{ D = "2.698";    E = "  +4.19%  " }
    46 = @{ D = "0.7405" }
    47 = @{ D = "4.219";    E = "  +0.43%  " }
    48 = @{ E = "  +2.56%  " }
    49 = @{ D = "1.003";    E = "  +0.02%  " }
    50 = @{ D = "0.08380";  E = "  +5.17%  " }
    51 = @{ D = "139.07";   E = "  +0.82%  " }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in @("B", "C", "D", "E")) {
        if ($cols.Contains($col)) {
            $cellRef = "$col$row"
            if ($col -eq "D") {
                Set-TextValue $cellRef $cols[$col]
            } else {
                $ws.Range($cellRef).Value = $cols[$col]
            }
        }
    }
}
